# Weekly update: insert a new price record for "Espinaca" (Vega Modelo de
# Temuco) as the new first row of the data block, pushing the existing
# rows 83-94 down to 84-95 (dimension grows from A1:R94 to A1:R95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83; everything currently at/after row 83 shifts down.
$ws.Rows("83:83").Insert()

$ws.Range("A83").Value2 = 10
$ws.Range("B83").Value2 = "Vega Modelo de Temuco"
$ws.Range("C83").Value2 = "La Araucanía"
$ws.Range("D83").Value2 = 44505
$ws.Range("E83").Value2 = 9
$ws.Range("F83").Value2 = 100112012
$ws.Range("G83").Value2 = "Espinaca"
$ws.Range("H83").Value2 = "Sin especificar"
$ws.Range("I83").Value2 = "Primera"
$ws.Range("J83").Value2 = 145
$ws.Range("K83").Value2 = 7000
$ws.Range("L83").Value2 = 8000
$ws.Range("M83").Value2 = 7448
$ws.Range("N83").Value2 = "$/docena de atados"
$ws.Range("O83").Value2 = "Región de La Araucanía"
$ws.Range("P83").Value2 = 2483
$ws.Range("Q83").Value2 = 3
$ws.Range("R83").Value2 = "Hortaliza"
